# Rename the four worksheets to include their units, as in the published
# version of the workbook.
$wb = $excel.ActiveWorkbook

$wsHeight  = $wb.Worksheets.Item(1)
$wsFourier = $wb.Worksheets.Item(2)
$wsCubic   = $wb.Worksheets.Item(3)
$wsTotal   = $wb.Worksheets.Item(4)

$wsHeight.Name  = "Height information (mm)"
$wsFourier.Name = "Fourier phase (rad)"
$wsCubic.Name   = "Cubic phase (rad)"
$wsTotal.Name   = "Total phase (rad)"

# Update the selected cell on the last (active) sheet to E8, matching the
# author's latest cursor position when the file was saved.
$wsTotal.Activate()
$wsTotal.Range("E8").Select()
